# AT01 Q5 Code of Ethics - body copy pass
# Adds spacing-after to several "MyStyle" paragraphs and fills in the
# previously-empty body paragraphs under each heading, plus a new
# (currently unused) "Me Signing" paragraph/character style pair.

$d = $word.ActiveDocument

function Set-ParagraphXml($para, [string]$innerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

# --- "Health & Safety" body paragraph: add spacing, tweak trailing runs ---
$p5 = '<w:p w14:paraId="2123F881" w14:textId="451E8920" w:rsidR="00413859" w:rsidRDefault="000603DD" w:rsidP="003661D1">' +
    '<w:pPr><w:pStyle w:val="MyStyle"/><w:spacing w:after="240"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This </w:t></w:r>' +
    '<w:r w:rsidR="006A5C2B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>policy shows the</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> commitment of </w:t></w:r>' +
    '<w:r w:rsidR="00F7217C"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Rainbow Hero Co' + [char]0x2019 + 's</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> management </w:t></w:r>' +
    '<w:r w:rsidR="00111351"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&amp;</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> workers to health &amp; safety</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. We aim to remove or reduce the risks to the health, safety &amp; welfare of all workers &amp; anyone else who may be affected by our business operations. </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>We also ensure all work activities are done as safely as possible.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d.Paragraphs(6) $p5

# --- "Honesty" body paragraph ---
$p7 = '<w:p w14:paraId="02884576" w14:textId="40DD65A2" w:rsidR="00413859" w:rsidRDefault="00413859" w:rsidP="00413859">' +
    '<w:pPr><w:pStyle w:val="MyStyle"/><w:spacing w:after="240"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Honesty is the best policy</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> &amp; we like to show that. To demonstrate honesty is by being transparent &amp; giving others permission to see the real you &amp; form their own opinion of you. Being </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>transparent means</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> you are predictable &amp; very obvious about your actions &amp; intentions.</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> To demonstrate honesty requires that you display a certain level of integrity.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d.Paragraphs(8) $p7

# --- "Privacy" body paragraph ---
$p9 = '<w:p w14:paraId="4660894B" w14:textId="007FEF59" w:rsidR="00413859" w:rsidRDefault="00413859" w:rsidP="00413859">' +
    '<w:pPr><w:pStyle w:val="MyStyle"/><w:spacing w:after="240"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>This policy is our ability to keep our personal information to ourselves and to control what happens if we share it with others</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. The only personal information </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>that we would want from our workers are names, addresses, phone number, emails address, &amp; tax file number. We won' + [char]0x2019 + 't share any of our workers personal information if they don' + [char]0x2019 + 't want it to be shared.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d.Paragraphs(10) $p9

# --- "Copyrights" body paragraph ---
$p11 = '<w:p w14:paraId="55AF82BE" w14:textId="4937A238" w:rsidR="00413859" w:rsidRDefault="00413859" w:rsidP="00413859">' +
    '<w:pPr><w:pStyle w:val="MyStyle"/><w:spacing w:after="240"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Materials used in connection with Rainbow Hero Co may be subject to copyright' + [char]0x00A9 + ' protection</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. The materials may include, but are not limited </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>to</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> documents, slides, images, audio, &amp; video.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d.Paragraphs(12) $p11

# --- "Professional Development/Personal Development" body paragraph (no lang rPr) ---
$p13 = '<w:p w14:paraId="60C2615E" w14:textId="1C1E3B73" w:rsidR="00413859" w:rsidRDefault="00413859" w:rsidP="00413859">' +
    '<w:pPr><w:pStyle w:val="MyStyle"/><w:spacing w:after="240"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">We encourage our workers to get more </w:t></w:r>' +
    '<w:r><w:t>professional development when possible.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $d.Paragraphs(14) $p13

# --- "Diversity & Inclusion": drop the extra blank paragraph, keep the last one ---
# (before: two trailing paragraphs after the heading - a bare blank para, then the
#  blank "MyStyle" para; after: just the single "MyStyle" para with spacing-after)
$d.Paragraphs(16).Range.Delete()
$d.Paragraphs(16).Range.ParagraphFormat.SpaceAfter = 12

# --- New (currently-unused) "Me Signing" paragraph style + its linked character style ---
$meSigning = $d.Styles.Add("Me Signing", 1)
$meSigning.QuickStyle = $true
$meSigning.Font.Name = "Brush Script MT"
$meSigning.Font.Color = 16711680
$meSigning.Font.LanguageID = "en-US"

$meSigningChar = $d.Styles.Add("Me Signing Char", 2)
$meSigningChar.BaseStyle = "DefaultParagraphFont"
$meSigningChar.Font.Name = "Brush Script MT"
$meSigningChar.Font.Color = 16711680
$meSigningChar.Font.LanguageID = "en-US"

$meSigning.LinkStyle = "Me Signing Char"
$meSigningChar.LinkStyle = "Me Signing"
